# SonoVerse workbook update
# 1) Fix capitalization of the clip-type label in row 3 (C3):
#    "Clip 1 B-mode + color + microV" -> "Clip 1 B-mode + Color + microV"
# 2) Insert a new "Pancreas" entry ("Chronic pancreatitis") as a new row,
#    keeping the existing alphabetical-by-Organ sort (inserted right
#    after the other Pancreas rows / before the Spleen rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix C3 capitalization ---
$ws.Range("C3").Value = "Clip 1 B-mode + Color + microV"

# --- Insert new row for "Chronic pancreatitis" at row 29 ---
$ws.Rows.Item(29).Insert()

$ws.Range("A29").Value = "Pancreas"
$ws.Range("B29").Value = "Chronic pancreatitis"
$ws.Range("C29").Value = "Clip 1 B-mode + Color"
$ws.Range("D29").Value = "https://youtu.be/VJdnjrAAO-4"
$ws.Hyperlinks.Add($ws.Range("D29"), "https://youtu.be/VJdnjrAAO-4")
$ws.Range("D29").Style = "Collegamento ipertestuale"

# --- Update selection to match the new working cell ---
$ws.Range("E21").Select()
